$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.13508
$ws.Range("H2").Value = 3.40524
$ws.Range("I2").Value = 0.1224366388308639
$ws.Range("J2").Value = 0.1224366388308639
$ws.Range("M2").Value = 2.781751333333333
$ws.Range("N2").Value = 8.345253999999999
$ws.Range("O2").Value = 0.06744008595411712
$ws.Range("P2").Value = 0.06744008595411713
$ws.Range("Q2").Value = 3.15751030344
$ws.Range("R2").Value = 28.41759273096
$ws.Range("S2").Value = 0.008257137446686658
$ws.Range("T2").Value = 0.008257137446686659

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.13508
$ws.Range("H3").Value = 3.40524
$ws.Range("I3").Value = 0.1224366388308639
$ws.Range("J3").Value = 0.1224366388308639
$ws.Range("O3").Value = 0.6170396927334101
$ws.Range("P3").Value = 0.6170396927334101
$ws.Range("Q3").Value = 28.88948256624
$ws.Range("R3").Value = 260.00534309616
$ws.Range("S3").Value = 0.07554826600350779
$ws.Range("T3").Value = 0.07554826600350779

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.13508
$ws.Range("H4").Value = 3.40524
$ws.Range("I4").Value = 0.1224366388308639
$ws.Range("J4").Value = 0.1224366388308639
$ws.Range("M4").Value = 3.163721333333334
$ws.Range("N4").Value = 9.491164000000001
$ws.Range("O4").Value = 0.07670047142539008
$ws.Range("P4").Value = 0.07670047142539008
$ws.Range("Q4").Value = 3.591076811040001
$ws.Range("R4").Value = 32.31969129936
$ws.Range("S4").Value = 0.009390947918067485
$ws.Range("T4").Value = 0.009390947918067485

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.13508
$ws.Range("H5").Value = 3.40524
$ws.Range("I5").Value = 0.1224366388308639
$ws.Range("J5").Value = 0.1224366388308639
$ws.Range("M5").Value = 4.241536
$ws.Range("N5").Value = 12.724608
$ws.Range("O5").Value = 0.1028307415511195
$ws.Range("P5").Value = 0.1028307415511195
$ws.Range("Q5").Value = 4.81448268288
$ws.Range("R5").Value = 43.33034414592
$ws.Range("S5").Value = 0.01259025036400434
$ws.Range("T5").Value = 0.01259025036400434

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.13508
$ws.Range("H6").Value = 3.40524
$ws.Range("I6").Value = 0.1224366388308639
$ws.Range("J6").Value = 0.1224366388308639
$ws.Range("M6").Value = 5.075855
$ws.Range("N6").Value = 15.227565
$ws.Range("O6").Value = 0.1230577634272013
$ws.Range("P6").Value = 0.1230577634272013
$ws.Range("Q6").Value = 5.7615014934
$ws.Range("R6").Value = 51.8535134406
$ws.Range("S6").Value = 0.01506677893607015
$ws.Range("T6").Value = 0.01506677893607015

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 1.13508
$ws.Range("H7").Value = 3.40524
$ws.Range("I7").Value = 0.1224366388308639
$ws.Range("J7").Value = 0.1224366388308639
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.5333846666666667
$ws.Range("N7").Value = 1.600154
$ws.Range("O7").Value = 0.01293124490876184
$ws.Range("P7").Value = 0.01293124490876184
$ws.Range("Q7").Value = 0.6054342674400002
$ws.Range("R7").Value = 5.44890840696
$ws.Range("S7").Value = 0.001583258162527521
$ws.Range("T7").Value = 0.001583258162527521

# Row 8
$ws.Range("G8").Value = 0.672624
$ws.Range("H8").Value = 2.017872
$ws.Range("I8").Value = 0.07255331937570129
$ws.Range("J8").Value = 0.07255331937570129
$ws.Range("M8").Value = 2.781751333333333
$ws.Range("N8").Value = 8.345253999999999
$ws.Range("O8").Value = 0.06744008595411712
$ws.Range("P8").Value = 0.06744008595411713
$ws.Range("Q8").Value = 1.871072708832
$ws.Range("R8").Value = 16.839654379488
$ws.Range("S8").Value = 0.004893002094953806
$ws.Range("T8").Value = 0.004893002094953807

# Row 9
$ws.Range("G9").Value = 0.672624
$ws.Range("H9").Value = 2.017872
$ws.Range("I9").Value = 0.07255331937570129
$ws.Range("J9").Value = 0.07255331937570129
$ws.Range("O9").Value = 0.6170396927334101
$ws.Range("P9").Value = 0.6170396927334101
$ws.Range("Q9").Value = 17.119286148672
$ws.Range("R9").Value = 154.073575338048
$ws.Range("S9").Value = 0.0447682778943717
$ws.Range("T9").Value = 0.0447682778943717

# Row 10
$ws.Range("G10").Value = 0.672624
$ws.Range("H10").Value = 2.017872
$ws.Range("I10").Value = 0.07255331937570129
$ws.Range("J10").Value = 0.07255331937570129
$ws.Range("M10").Value = 3.163721333333334
$ws.Range("N10").Value = 9.491164000000001
$ws.Range("O10").Value = 0.07670047142539008
$ws.Range("P10").Value = 0.07670047142539008
$ws.Range("Q10").Value = 2.127994898112
$ws.Range("R10").Value = 19.151954083008
$ws.Range("S10").Value = 0.005564873799593177
$ws.Range("T10").Value = 0.005564873799593177

# Row 11
$ws.Range("G11").Value = 0.672624
$ws.Range("H11").Value = 2.017872
$ws.Range("I11").Value = 0.07255331937570129
$ws.Range("J11").Value = 0.07255331937570129
$ws.Range("M11").Value = 4.241536
$ws.Range("N11").Value = 12.724608
$ws.Range("O11").Value = 0.1028307415511195
$ws.Range("P11").Value = 0.1028307415511195
$ws.Range("Q11").Value = 2.852958910464
$ws.Range("R11").Value = 25.676630194176
$ws.Range("S11").Value = 0.007460711633398573
$ws.Range("T11").Value = 0.007460711633398573

# Row 12
$ws.Range("G12").Value = 0.672624
$ws.Range("H12").Value = 2.017872
$ws.Range("I12").Value = 0.07255331937570129
$ws.Range("J12").Value = 0.07255331937570129
$ws.Range("M12").Value = 5.075855
$ws.Range("N12").Value = 15.227565
$ws.Range("O12").Value = 0.1230577634272013
$ws.Range("P12").Value = 0.1230577634272013
$ws.Range("Q12").Value = 3.41414189352
$ws.Range("R12").Value = 30.72727704168
$ws.Range("S12").Value = 0.008928249211593232
$ws.Range("T12").Value = 0.008928249211593233

# Row 13
$ws.Range("G13").Value = 0.672624
$ws.Range("H13").Value = 2.017872
$ws.Range("I13").Value = 0.07255331937570129
$ws.Range("J13").Value = 0.07255331937570129
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.5333846666666667
$ws.Range("N13").Value = 1.600154
$ws.Range("O13").Value = 0.01293124490876184
$ws.Range("P13").Value = 0.01293124490876184
$ws.Range("Q13").Value = 0.3587673280320001
$ws.Range("R13").Value = 3.228905952288
$ws.Range("S13").Value = 0.000938204741790809
$ws.Range("T13").Value = 0.0009382047417908089

# Row 14
$ws.Range("G14").Value = 7.463050333333332
$ws.Range("H14").Value = 22.389151
$ws.Range("I14").Value = 0.8050100417934347
$ws.Range("J14").Value = 0.8050100417934348
$ws.Range("M14").Value = 2.781751333333333
$ws.Range("N14").Value = 8.345253999999999
$ws.Range("O14").Value = 0.06744008595411712
$ws.Range("P14").Value = 0.06744008595411713
$ws.Range("Q14").Value = 20.76035021548377
$ws.Range("R14").Value = 186.843151939354
$ws.Range("S14").Value = 0.05428994641247665
$ws.Range("T14").Value = 0.05428994641247667

# Row 15
$ws.Range("G15").Value = 7.463050333333332
$ws.Range("H15").Value = 22.389151
$ws.Range("I15").Value = 0.8050100417934347
$ws.Range("J15").Value = 0.8050100417934348
$ws.Range("O15").Value = 0.6170396927334101
$ws.Range("P15").Value = 0.6170396927334101
$ws.Range("Q15").Value = 189.9457857558982
$ws.Range("R15").Value = 1709.512071803084
$ws.Range("S15").Value = 0.4967231488355306
$ws.Range("T15").Value = 0.4967231488355306

# Row 16
$ws.Range("G16").Value = 7.463050333333332
$ws.Range("H16").Value = 22.389151
$ws.Range("I16").Value = 0.8050100417934347
$ws.Range("J16").Value = 0.8050100417934348
$ws.Range("M16").Value = 3.163721333333334
$ws.Range("N16").Value = 9.491164000000001
$ws.Range("O16").Value = 0.07670047142539008
$ws.Range("P16").Value = 0.07670047142539008
$ws.Range("Q16").Value = 23.61101155130711
$ws.Range("R16").Value = 212.499103961764
$ws.Range("S16").Value = 0.06174464970772941
$ws.Range("T16").Value = 0.06174464970772942

# Row 17
$ws.Range("G17").Value = 7.463050333333332
$ws.Range("H17").Value = 22.389151
$ws.Range("I17").Value = 0.8050100417934347
$ws.Range("J17").Value = 0.8050100417934348
$ws.Range("M17").Value = 4.241536
$ws.Range("N17").Value = 12.724608
$ws.Range("O17").Value = 0.1028307415511195
$ws.Range("P17").Value = 0.1028307415511195
$ws.Range("Q17").Value = 31.65479665864533
$ws.Range("R17").Value = 284.893169927808
$ws.Range("S17").Value = 0.08277977955371663
$ws.Range("T17").Value = 0.08277977955371664

# Row 18
$ws.Range("G18").Value = 7.463050333333332
$ws.Range("H18").Value = 22.389151
$ws.Range("I18").Value = 0.8050100417934347
$ws.Range("J18").Value = 0.8050100417934348
$ws.Range("M18").Value = 5.075855
$ws.Range("N18").Value = 15.227565
$ws.Range("O18").Value = 0.1230577634272013
$ws.Range("P18").Value = 0.1230577634272013
$ws.Range("Q18").Value = 37.88136134970166
$ws.Range("R18").Value = 340.932252147315
$ws.Range("S18").Value = 0.09906273527953795
$ws.Range("T18").Value = 0.09906273527953798

# Row 19
$ws.Range("G19").Value = 7.463050333333332
$ws.Range("H19").Value = 22.389151
$ws.Range("I19").Value = 0.8050100417934347
$ws.Range("J19").Value = 0.8050100417934348
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.5333846666666667
$ws.Range("N19").Value = 1.600154
$ws.Range("O19").Value = 0.01293124490876184
$ws.Range("P19").Value = 0.01293124490876184
$ws.Range("Q19").Value = 3.980676614361556
$ws.Range("R19").Value = 35.826089529254
$ws.Range("S19").Value = 0.01040978200444351
$ws.Range("T19").Value = 0.01040978200444351
